$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("EXERCISE 1")
$ws2 = $wb.Worksheets.Item("EXERCISE 2")

# Insert a new column before column B ("COMMENTI" will be filled in later so
# that it lands last in the shared-string table, matching the order strings
# were authored in Excel).
$ws1.Columns("B").Insert()
$ws1.Columns("B").ColumnWidth = 29.833333333333332

# New "LANCIO n" columns appended after the existing data (now K is the last
# used column after the insert above).
$ws1.Range("L1").Value = "LANCIO 7"
$ws1.Range("M1").Value = "LANCIO 8"
$ws1.Range("N1").Value = "LANCIO 9 "
$ws1.Range("O1").Value = "LANCIO 10"

$ws1.Columns("L").ColumnWidth = 17
$ws1.Columns("M").ColumnWidth = 12.333333333333334
$ws1.Columns("N").ColumnWidth = 17.833333333333332
$ws1.Columns("O").ColumnWidth = 15.5

# New row with the "Hierarchical OPT (GIGI)" results.
$ws1.Range("A3").Value = "Hierarchical OPT (GIGI)"
$ws1.Range("F3").Value = 4.2328977584838796
$ws1.Range("G3").Value = 4.3476593494415203
$ws1.Range("H3").Value = 4.2484049797058097
$ws1.Range("I3").Value = 4.2112441062927202
$ws1.Range("J3").Value = 4.18296027183532
$ws1.Range("K3").Value = 4.5008549690246502
$ws1.Range("L3").Value = 4.2073841094970703
$ws1.Range("M3").Value = 4.4207534790039
$ws1.Range("N3").Value = 4.6460390090942303
$ws1.Range("O3").Value = 4.1695647239684996
$ws1.Range("E3").Formula = "=AVERAGE(F3:O3)"

# "COMMENTI" column header - added last so it becomes the final shared
# string in the table.
$ws1.Range("B1").Value = "COMMENTI"

# Selection / active sheet bookkeeping: EXERCISE 1 becomes the tab in focus,
# EXERCISE 2 loses its tabSelected flag, selection moves to E3.
[void]$ws2.Range("A12").Select()
[void]$ws1.Activate()
[void]$ws1.Range("E3").Select()
